$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Appointment AP002 (row 3) has had its outcome recorded: mark its
# "Appointment Status" (column E) as Done instead of Pending.
$ws.Range("E3").Value = "Done"
